$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D), Volume(1h) (E), and Hora (G) columns for rows 2-51
# as scraped by the GitHub Actions symbol-list update job.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '274.68'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.98%'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '9'

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.49'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-2.81%'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '9'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.881'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '2.02%'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '9'

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06335'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.86%'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '9'

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.887'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.66%'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '9'

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.310'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.16%'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '9'

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.259'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '33.29%'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '9'

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8708'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.76%'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '9'

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1522'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '4.64%'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '9'

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.05031'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-2.73%'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '9'

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07400'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '1.57%'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '9'

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03025'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-2.64%'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '9'

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09031'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.39%'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '9'

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.18%'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '9'

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0006330'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.74%'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '9'

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.005877'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.75%'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '9'

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.450'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.10%'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '9'

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.54%'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '9'

# Row 20
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '9'

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1324'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.85%'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '9'

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.914'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.59%'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '9'

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04361'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.75%'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '9'

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001180'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.22%'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '9'

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004249'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.63%'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '9'

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.24%'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '9'

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001677'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-0.89%'
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '9'

# Row 28
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '9'

# Row 29
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '9'

# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '9'

# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '9'

# Row 32
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '9'

# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '9'

# Row 34
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '9'

# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '9'

# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '9'

# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '9'

# Row 38
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '9'

# Row 39
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '9'

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04102'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.93%'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '9'

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006967'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '8.74%'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '9'

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1167'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '1.06%'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '9'

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-2.68%'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '9'

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-9.41%'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '9'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005262'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '2.92%'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '9'

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.486'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-37.49%'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '9'

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.01998'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-11.32%'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '9'

# Row 48
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '9'

# Row 49
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '9'

# Row 50
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '9'

# Row 51
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '9'
